$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions): update "想去人数" (want-to-go count) column F
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F4").Value = 1640   # 南宁·草莓动漫节: 1635 -> 1640
$wsExhibit.Range("F6").Value = 57     # 南宁·布谷鸟动漫展4th: 56 -> 57

# Sheet "全部类型" (all types): same events appear with different row offsets
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value = 1640       # 南宁·草莓动漫节: 1635 -> 1640
$wsAll.Range("F7").Value = 57         # 南宁·布谷鸟动漫展4th: 56 -> 57
